$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Hex") to hold the Background_Color data.
# This shifts Money, Min_Money, Peek, Min_Peek, Trade, Min_Trade one column to the right.
$ws.Columns.Item(3).Insert()

# New header
$ws.Range("C1").Value = "Background_Color"

# New light/washed-out background color values for each row (lighter version of the Hex color)
$ws.Range("C2").Value = "f6c8c8"
$ws.Range("C3").Value = "fff8a1"
$ws.Range("C4").Value = "d1ffb7"
$ws.Range("C5").Value = "bbc3fb"
$ws.Range("C6").Value = "b5b2af"

# Update the active selection to match the edited cell
$ws.Range("C2").Select()
